$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 270, shifting the existing rows 270-280 down to 271-281
$ws.Rows(270).Insert()

# Populate the newly inserted row 270 with the new weekly record
$ws.Cells.Item(270, 1).Value = 9
$ws.Cells.Item(270, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(270, 3).Value = "Metropolitana"
$ws.Cells.Item(270, 4).Value = 44509
$ws.Cells.Item(270, 5).Value = 13
$ws.Cells.Item(270, 6).Value = 100112012
$ws.Cells.Item(270, 7).Value = "Espinaca"
$ws.Cells.Item(270, 8).Value = "Sin especificar"
$ws.Cells.Item(270, 9).Value = "Primera"
$ws.Cells.Item(270, 10).Value = 250
$ws.Cells.Item(270, 11).Value = 6000
$ws.Cells.Item(270, 12).Value = 7000
$ws.Cells.Item(270, 13).Value = 6500
$ws.Cells.Item(270, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(270, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(270, 16).Value = 650
$ws.Cells.Item(270, 17).Value = 10
$ws.Cells.Item(270, 18).Value = "Hortaliza"
